$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 449, shifting existing rows 449..561 down to 450..562
$ws.Rows.Item(449).Insert()

# Populate the newly inserted row 449 with the new data record
$ws.Cells.Item(449, 1).Value = 6
$ws.Cells.Item(449, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(449, 3).Value = "Metropolitana"
$ws.Cells.Item(449, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 9 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(449, 5).Value = 13
$ws.Cells.Item(449, 6).Value = 100112044
$ws.Cells.Item(449, 7).Value = "Perejil"
$ws.Cells.Item(449, 8).Value = "Sin especificar"
$ws.Cells.Item(449, 9).Value = "Primera"
$ws.Cells.Item(449, 10).Value = 220
$ws.Cells.Item(449, 11).Value = 18000
$ws.Cells.Item(449, 12).Value = 20000
$ws.Cells.Item(449, 13).Value = 19091
$ws.Cells.Item(449, 14).Value = "`$/docena de atados"
$ws.Cells.Item(449, 15).Value = "Región Metropolitana"
$ws.Cells.Item(449, 16).Value = 6364
$ws.Cells.Item(449, 17).Value = 3
$ws.Cells.Item(449, 18).Value = "Hortaliza"
